# Auto-update the water delivery tracker sheet:
# For each data row, decrement the "remaining days" (column E) by 1.
# When remaining days would reach 0, reset it to the "total days" value
# (column D) and bump the "start date" (column F) to the current cycle
# start date (2026-01-20). Rows whose start date is not a clean 8-digit
# date (e.g. malformed data) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; E = 6; F = 20260112 },
    @{ Row = 3; E = 6; F = 20260112 },
    @{ Row = 4; E = 6; F = 20260112 },
    @{ Row = 5; E = 4; F = 20260114 },
    @{ Row = 6; E = 6; F = 20260112 },
    @{ Row = 7; E = 4; F = 20260114 },
    @{ Row = 8; E = 6; F = 20260112 },
    @{ Row = 9; E = 4; F = 20260114 },
    @{ Row = 10; E = 6; F = 20260119 },
    @{ Row = 11; E = 6; F = 20260112 },
    @{ Row = 12; E = 4; F = 20260114 },
    @{ Row = 13; E = 6; F = 20260112 },
    @{ Row = 14; E = 6; F = 20260112 },
    @{ Row = 15; E = 6; F = 20260112 },
    @{ Row = 16; E = 8; F = 20260118 },
    @{ Row = 17; E = 4; F = 20260114 },
    @{ Row = 18; E = 7; F = 20260117 },
    @{ Row = 19; E = 7; F = 20260117 },
    @{ Row = 20; E = 7; F = 20260117 },
    @{ Row = 21; E = 7; F = 20260117 },
    @{ Row = 22; E = 4; F = 20260114 },
    @{ Row = 23; E = 4; F = 20260114 },
    @{ Row = 24; E = 4; F = 20260114 },
    @{ Row = 25; E = 4; F = 20260114 },
    @{ Row = 26; E = 4; F = 20260114 },
    @{ Row = 27; E = 7; F = 20260120 },
    @{ Row = 28; E = 7; F = 20260117 },
    @{ Row = 29; E = 7; F = 20260117 },
    @{ Row = 30; E = 7; F = 20260117 },
    @{ Row = 31; E = 7; F = 20260117 },
    @{ Row = 32; E = 7; F = 20260117 },
    @{ Row = 33; E = 7; F = 20260117 },
    @{ Row = 34; E = 7; F = 20260117 },
    @{ Row = 35; E = 7; F = 20260117 },
    @{ Row = 37; E = 7; F = 20260117 },
    @{ Row = 38; E = 7; F = 20260117 },
    @{ Row = 39; E = 7; F = 20260117 },
    @{ Row = 40; E = 6; F = 20260119 },
    @{ Row = 41; E = 6; F = 20260119 },
    @{ Row = 42; E = 7; F = 20260117 },
    @{ Row = 43; E = 4; F = 20260114 },
    @{ Row = 44; E = 6; F = 20260119 },
    @{ Row = 45; E = 4; F = 20260114 },
    @{ Row = 46; E = 6; F = 20260119 },
    @{ Row = 47; E = 7; F = 20260117 },
    @{ Row = 48; E = 6; F = 20260119 },
    @{ Row = 49; E = 7; F = 20260120 },
    @{ Row = 50; E = 2; F = 20260112 },
    @{ Row = 51; E = 2; F = 20260112 },
    @{ Row = 52; E = 2; F = 20260112 },
    @{ Row = 53; E = 2; F = 20260112 },
    @{ Row = 54; E = 2; F = 20260112 },
    @{ Row = 55; E = 2; F = 20260112 },
    @{ Row = 56; E = 2; F = 20260112 },
    @{ Row = 57; E = 2; F = 20260112 },
    @{ Row = 58; E = 6; F = 20260116 },
    @{ Row = 59; E = 6; F = 20260116 },
    @{ Row = 60; E = 6; F = 20260116 },
    @{ Row = 61; E = 7; F = 20260120 },
    @{ Row = 62; E = 6; F = 20260116 },
    @{ Row = 63; E = 6; F = 20260116 },
    @{ Row = 64; E = 6; F = 20260116 },
    @{ Row = 65; E = 7; F = 20260117 },
    @{ Row = 66; E = 7; F = 20260117 },
    @{ Row = 67; E = 7; F = 20260117 },
    @{ Row = 68; E = 7; F = 20260117 },
    @{ Row = 69; E = 7; F = 20260117 },
    @{ Row = 70; E = 8; F = 20260118 },
    @{ Row = 71; E = 8; F = 20260118 },
    @{ Row = 72; E = 8; F = 20260118 },
    @{ Row = 73; E = 8; F = 20260118 },
    @{ Row = 74; E = 8; F = 20260118 },
    @{ Row = 75; E = 8; F = 20260118 },
    @{ Row = 76; E = 8; F = 20260118 },
    @{ Row = 77; E = 1; F = 20260111 },
    @{ Row = 78; E = 1; F = 20260111 },
    @{ Row = 79; E = 1; F = 20260111 },
    @{ Row = 80; E = 1; F = 20260111 },
    @{ Row = 81; E = 1; F = 20260111 },
    @{ Row = 82; E = 1; F = 20260111 },
    @{ Row = 83; E = 1; F = 20260111 },
    @{ Row = 84; E = 1; F = 20260111 },
    @{ Row = 85; E = 1; F = 20260111 },
    @{ Row = 86; E = 1; F = 20260111 },
    @{ Row = 87; E = 6; F = 20260119 },
    @{ Row = 88; E = 6; F = 20260119 },
    @{ Row = 89; E = 6; F = 20260119 },
    @{ Row = 90; E = 6; F = 20260119 },
    @{ Row = 91; E = 4; F = 20260114 },
    @{ Row = 92; E = 6; F = 20260119 },
    @{ Row = 93; E = 1; F = 20260111 },
    @{ Row = 94; E = 2; F = 20260115 },
    @{ Row = 95; E = 10; F = 20260120 },
    @{ Row = 96; E = 8; F = 20260118 },
    @{ Row = 97; E = 8; F = 20260118 },
    @{ Row = 98; E = 8; F = 20260118 },
    @{ Row = 99; E = 8; F = 20260118 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
